$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "edit1"
$ws.Range("B14").Value = "riya-morankar"
$ws.Range("C14").Value = "Squashed"
$ws.Range("D14").Value = "N/A"

# "2025-06-18" looks like a date, so Excel would normally auto-convert it
# to a date serial number with a date NumberFormat when assigned directly.
# Force it to be stored as plain text (matching the other "Date" column
# cells in this sheet, which are inline/text strings, not real dates):
# temporarily mark the cell as Text, assign the literal string, then
# clear the formatting we introduced so no stray style lingers on the cell.
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2025-06-18"
$ws.Range("E14").ClearFormats()

$ws.Range("F14").Value = "2efdb7acf2306cf36253f988e5a0ebd5022b326d"
